$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'57.839.00"
$ws.Range("E2").Value = "  +2.68%  "

# Row 3
$ws.Range("D3").Value = "'3.043.31"
$ws.Range("E3").Value = "  +2.28%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'516.83"
$ws.Range("E5").Value = "  +2.68%  "

# Row 6
$ws.Range("D6").Value = "'141.10"
$ws.Range("E6").Value = "  +4.76%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.444"
$ws.Range("E8").Value = "  +3.70%  "

# Row 9
$ws.Range("E9").Value = "  +2.56%  "

# Row 10
$ws.Range("D10").Value = "'0.111"
$ws.Range("E10").Value = "  +4.45%  "

# Row 11
$ws.Range("E11").Value = "  +5.34%  "

# Row 12
$ws.Range("D12").Value = "'3.572.46"
$ws.Range("E12").Value = "  +2.47%  "

# Row 13
$ws.Range("E13").Value = "  +2.33%  "

# Row 14
$ws.Range("D14").Value = "'26.78"
$ws.Range("E14").Value = "  +6.48%  "

# Row 15
$ws.Range("D15").Value = "'0.0000168"
$ws.Range("E15").Value = "  +11.02%  "

# Row 16
$ws.Range("D16").Value = "'57.871.79"
$ws.Range("E16").Value = "  +2.74%  "

# Row 17
$ws.Range("E17").Value = "  +10.01%  "

# Row 18
$ws.Range("D18").Value = "'3.046.48"
$ws.Range("E18").Value = "  +2.33%  "

# Row 19
$ws.Range("D19").Value = "'13.00"
$ws.Range("E19").Value = "  +5.37%  "

# Row 20
$ws.Range("D20").Value = "'8.06"
$ws.Range("E20").Value = "  +3.89%  "

# Row 21
$ws.Range("D21").Value = "'337.05"
$ws.Range("E21").Value = "  +4.06%  "

# Row 22
$ws.Range("D22").Value = "'5.77"
$ws.Range("E22").Value = "  +1.40%  "

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24
$ws.Range("D24").Value = "'0.501"
$ws.Range("E24").Value = "  +6.54%  "

# Row 25
$ws.Range("D25").Value = "'64.98"
$ws.Range("E25").Value = "  +4.97%  "

# Row 26
$ws.Range("E26").Value = "  +3.94%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.34%  "

# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0₃0944"
$ws.Range("E28").Value = "  +6.75%  "

# Row 29
$ws.Range("D29").Value = "'6.88"
$ws.Range("E29").Value = "  +6.23%  "

# Row 30
$ws.Range("D30").Value = "'7.53"
$ws.Range("E30").Value = "  +10.98%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.82"
$ws.Range("E31").Value = "  +4.30%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.23"
$ws.Range("E32").Value = "  +3.05%  "

# Row 33
$ws.Range("D33").Value = "'20.97"
$ws.Range("E33").Value = "  +2.40%  "

# Row 34
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'156.96"
$ws.Range("E34").Value = "  -0.79%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.75"
$ws.Range("E35").Value = "  +6.65%  "

# Row 36
$ws.Range("D36").Value = "'5.90"
$ws.Range("E36").Value = "  +6.43%  "

# Row 37
$ws.Range("D37").Value = "'1.29"
$ws.Range("E37").Value = "  +2.01%  "

# Row 38
$ws.Range("D38").Value = "'25.12"
$ws.Range("E38").Value = "  +9.69%  "

# Row 39
$ws.Range("D39").Value = "'0.0690"
$ws.Range("E39").Value = "  +2.57%  "

# Row 40
$ws.Range("D40").Value = "'3.080.22"
$ws.Range("E40").Value = "  +2.36%  "

# Row 41
$ws.Range("D41").Value = "'37.82"
$ws.Range("E41").Value = "  +4.46%  "

# Row 42
$ws.Range("D42").Value = "'3.89"
$ws.Range("E42").Value = "  +9.56%  "

# Row 43
$ws.Range("E43").Value = "  +0.10%  "

# Row 44
$ws.Range("D44").Value = "'0.663"
$ws.Range("E44").Value = "  +3.61%  "

# Row 45
$ws.Range("D45").Value = "'2.311.58"
$ws.Range("E45").Value = "  +2.97%  "

# Row 46
$ws.Range("D46").Value = "'1.44"
$ws.Range("E46").Value = "  +3.53%  "

# Row 47
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +2.02%  "

# Row 48
$ws.Range("D48").Value = "'6.07"
$ws.Range("E48").Value = "  +5.11%  "

# Row 49
$ws.Range("D49").Value = "'0.0241"
$ws.Range("E49").Value = "  +2.79%  "

# Row 50
$ws.Range("D50").Value = "'19.68"
$ws.Range("E50").Value = "  +3.92%  "

# Row 51
$ws.Range("D51").Value = "'1.85"
$ws.Range("E51").Value = "  -4.40%  "

